$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2210 = "_FV2210"
$fv2304 = "_FV2304"

$headers = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i] + $fv2210
}

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $headers[$i] + $fv2304
}

$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U82"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

